$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the date-like cells we are about to write so
# Excel keeps them as literal text instead of converting to serial date
# numbers (only the cells whose values actually change/are new).
$ws.Range("N3:O3").NumberFormat = "@"
$ws.Range("J4:K4").NumberFormat = "@"
$ws.Range("N4:O4").NumberFormat = "@"

$ws.Range("A3").Value = "LV-1753526987168-INVG"
$ws.Range("B3").Value = "Manager"
$ws.Range("G3").Value = "HR Executive "
$ws.Range("H3").Value = "Human Resource"
$ws.Range("N3").Value = "2025-07-26"
$ws.Range("O3").Value = "2025-07-26"

$ws.Range("A4").Value = "LV-1753528055049-S2BW"
$ws.Range("B4").Value = "Manager"
$ws.Range("C4").Value = "PILLP305"
$ws.Range("D4").Value = "Animesh Roy"
$ws.Range("E4").Value = 8145312848
$ws.Range("F4").Value = "animesh.roy@pillp.in"
$ws.Range("G4").Value = "Mechanical HOD"
$ws.Range("H4").Value = "Mechanical"
$ws.Range("I4").Value = "Raichur"
$ws.Range("J4").Value = "02-12-2025"
$ws.Range("K4").Value = "07-12-2025"
$ws.Range("L4").Value = "For my personal reason"
$ws.Range("M4").Value = "Pending"
$ws.Range("N4").Value = "2025-07-26"
$ws.Range("O4").Value = "2025-07-26"
